# Updates the weekly FlashScore odds sheet for 2024-10-10:
#  - Row 2 becomes a newly-scraped match (Id 65BJiW3D, Ind. Medellin vs Alianza).
#  - The match that used to sit in row 2 (Santa Fe vs Chico) moves to row 3, with refreshed odds.
#  - The match that used to sit in row 3 (Correcaminos vs Atl. Morelia) moves to row 4, with refreshed odds.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 so the previous rows 2 and 3 shift down to 3 and 4.
$ws.Rows.Item(3).Insert()

# Row 2: new match -- Id 65BJiW3D, Ind. Medellin vs Alianza.
$ws.Cells.Item(2, 1).Value = '65BJiW3D'  # A2
$ws.Cells.Item(2, 2).Value = '''10/10/2024'  # B2
$ws.Cells.Item(2, 3).Value = '20:30'  # C2
$ws.Cells.Item(2, 4).Value = 'COLOMBIA - PRIMERA A'  # D2
$ws.Cells.Item(2, 5).Value = 'Ind. Medellin'  # E2
$ws.Cells.Item(2, 6).Value = 'Alianza'  # F2
$ws.Cells.Item(2, 7).Value = 1.65  # G2
$ws.Cells.Item(2, 8).Value = 3.8  # H2
$ws.Cells.Item(2, 9).Value = 5.25  # I2
$ws.Cells.Item(2, 10).Value = 2.3  # J2
$ws.Cells.Item(2, 11).Value = 2.1  # K2
$ws.Cells.Item(2, 12).Value = 6  # L2
$ws.Cells.Item(2, 13).Value = 1.07  # M2
$ws.Cells.Item(2, 14).Value = 9  # N2
$ws.Cells.Item(2, 15).Value = 1.4  # O2
$ws.Cells.Item(2, 16).Value = 2.75  # P2
$ws.Cells.Item(2, 17).Value = 2.2  # Q2
$ws.Cells.Item(2, 18).Value = 1.65  # R2
$ws.Cells.Item(2, 19).Value = 1.44  # S2
$ws.Cells.Item(2, 20).Value = 2.63  # T2
$ws.Cells.Item(2, 21).Value = 2.2  # U2
$ws.Cells.Item(2, 22).Value = 1.62  # V2
$ws.Cells.Item(2, 23).Value = 5.5  # W2
$ws.Cells.Item(2, 24).Value = 7  # X2
$ws.Cells.Item(2, 25).Value = 9  # Y2
$ws.Cells.Item(2, 26).Value = 12  # Z2
$ws.Cells.Item(2, 27).Value = 15  # AA2
$ws.Cells.Item(2, 28).Value = 34  # AB2
$ws.Cells.Item(2, 29).Value = 8  # AC2
$ws.Cells.Item(2, 30).Value = 7  # AD2
$ws.Cells.Item(2, 31).Value = 21  # AE2
$ws.Cells.Item(2, 32).Value = 81  # AF2
$ws.Cells.Item(2, 33).Value = 900  # AG2
$ws.Cells.Item(2, 34).Value = 11  # AH2
$ws.Cells.Item(2, 35).Value = 26  # AI2
$ws.Cells.Item(2, 36).Value = 17  # AJ2
$ws.Cells.Item(2, 37).Value = 51  # AK2
$ws.Cells.Item(2, 38).Value = 41  # AL2
$ws.Cells.Item(2, 39).Value = 51  # AM2
$ws.Cells.Item(2, 40).Value = 3.5  # AN2
$ws.Cells.Item(2, 41).Value = 9  # AO2
$ws.Cells.Item(2, 42).Value = 23  # AP2
$ws.Cells.Item(2, 43).Value = 29  # AQ2
$ws.Cells.Item(2, 44).Value = 51  # AR2
$ws.Cells.Item(2, 45).Value = 201  # AS2
$ws.Cells.Item(2, 46).Value = 2.63  # AT2
$ws.Cells.Item(2, 47).Value = 9.5  # AU2
$ws.Cells.Item(2, 48).Value = 67  # AV2
$ws.Cells.Item(2, 49).Value = 7  # AW2
$ws.Cells.Item(2, 50).Value = 34  # AX2
$ws.Cells.Item(2, 51).Value = 41  # AY2
$ws.Cells.Item(2, 52).Value = 126  # AZ2
$ws.Cells.Item(2, 53).Value = 151  # BA2
$ws.Cells.Item(2, 54).Value = 351  # BB2
$ws.Cells.Item(2, 55).Value = 126  # BC2
$ws.Cells.Item(2, 56).Value = 126  # BD2

# Row 3: Santa Fe vs Chico (Id 8UHmSCuQ), shifted down from row 2, odds refreshed.
$ws.Cells.Item(3, 1).Value = '8UHmSCuQ'  # A3
$ws.Cells.Item(3, 2).Value = '''10/10/2024'  # B3
$ws.Cells.Item(3, 3).Value = '22:30'  # C3
$ws.Cells.Item(3, 4).Value = 'COLOMBIA - PRIMERA A'  # D3
$ws.Cells.Item(3, 5).Value = 'Santa Fe'  # E3
$ws.Cells.Item(3, 6).Value = 'Chico'  # F3
$ws.Cells.Item(3, 7).Value = 1.44  # G3
$ws.Cells.Item(3, 8).Value = 4.1  # H3
$ws.Cells.Item(3, 9).Value = 7.5  # I3
$ws.Cells.Item(3, 10).Value = 2.05  # J3
$ws.Cells.Item(3, 11).Value = 2.1  # K3
$ws.Cells.Item(3, 12).Value = 8  # L3
$ws.Cells.Item(3, 13).Value = 1.07  # M3
$ws.Cells.Item(3, 14).Value = 9  # N3
$ws.Cells.Item(3, 15).Value = 1.4  # O3
$ws.Cells.Item(3, 16).Value = 2.75  # P3
$ws.Cells.Item(3, 17).Value = 2.2  # Q3
$ws.Cells.Item(3, 18).Value = 1.65  # R3
$ws.Cells.Item(3, 19).Value = 1.5  # S3
$ws.Cells.Item(3, 20).Value = 2.5  # T3
$ws.Cells.Item(3, 21).Value = 2.5  # U3
$ws.Cells.Item(3, 22).Value = 1.5  # V3
$ws.Cells.Item(3, 23).Value = 5  # W3
$ws.Cells.Item(3, 24).Value = 5.5  # X3
$ws.Cells.Item(3, 25).Value = 9.5  # Y3
$ws.Cells.Item(3, 26).Value = 9  # Z3
$ws.Cells.Item(3, 27).Value = 15  # AA3
$ws.Cells.Item(3, 28).Value = 41  # AB3
$ws.Cells.Item(3, 29).Value = 8  # AC3
$ws.Cells.Item(3, 30).Value = 8.5  # AD3
$ws.Cells.Item(3, 31).Value = 29  # AE3
$ws.Cells.Item(3, 32).Value = 101  # AF3
$ws.Cells.Item(3, 33).Value = 201  # AG3
$ws.Cells.Item(3, 34).Value = 13  # AH3
$ws.Cells.Item(3, 35).Value = 34  # AI3
$ws.Cells.Item(3, 36).Value = 23  # AJ3
$ws.Cells.Item(3, 37).Value = 101  # AK3
$ws.Cells.Item(3, 38).Value = 67  # AL3
$ws.Cells.Item(3, 39).Value = 67  # AM3
$ws.Cells.Item(3, 40).Value = 3.2  # AN3
$ws.Cells.Item(3, 41).Value = 7.5  # AO3
$ws.Cells.Item(3, 42).Value = 23  # AP3
$ws.Cells.Item(3, 43).Value = 23  # AQ3
$ws.Cells.Item(3, 44).Value = 51  # AR3
$ws.Cells.Item(3, 45).Value = 251  # AS3
$ws.Cells.Item(3, 46).Value = 2.5  # AT3
$ws.Cells.Item(3, 47).Value = 11  # AU3
$ws.Cells.Item(3, 48).Value = 81  # AV3
$ws.Cells.Item(3, 49).Value = 8.5  # AW3
$ws.Cells.Item(3, 50).Value = 41  # AX3
$ws.Cells.Item(3, 51).Value = 51  # AY3
$ws.Cells.Item(3, 52).Value = 201  # AZ3
$ws.Cells.Item(3, 53).Value = 251  # BA3
$ws.Cells.Item(3, 54).Value = 501  # BB3
$ws.Cells.Item(3, 55).Value = 126  # BC3
$ws.Cells.Item(3, 56).Value = 126  # BD3

# Row 4: Correcaminos vs Atl. Morelia (Id hCptA7hl), shifted down from row 3, odds refreshed.
$ws.Cells.Item(4, 1).Value = 'hCptA7hl'  # A4
$ws.Cells.Item(4, 2).Value = '''10/10/2024'  # B4
$ws.Cells.Item(4, 3).Value = '22:00'  # C4
$ws.Cells.Item(4, 4).Value = 'MEXICO - LIGA DE EXPANSION MX'  # D4
$ws.Cells.Item(4, 5).Value = 'Correcaminos'  # E4
$ws.Cells.Item(4, 6).Value = 'Atl. Morelia'  # F4
$ws.Cells.Item(4, 7).Value = 2.77  # G4
$ws.Cells.Item(4, 8).Value = 3.05  # H4
$ws.Cells.Item(4, 9).Value = 2.47  # I4
$ws.Cells.Item(4, 10).Value = 3.25  # J4
$ws.Cells.Item(4, 11).Value = 2.07  # K4
$ws.Cells.Item(4, 12).Value = 3.05  # L4
$ws.Cells.Item(4, 13).Value = 1.01  # M4
$ws.Cells.Item(4, 14).Value = 8  # N4
$ws.Cells.Item(4, 15).Value = 1.31  # O4
$ws.Cells.Item(4, 16).Value = 2.9  # P4
$ws.Cells.Item(4, 17).Value = 1.98  # Q4
$ws.Cells.Item(4, 18).Value = 1.75  # R4
$ws.Cells.Item(4, 19).Value = 1.42  # S4
$ws.Cells.Item(4, 20).Value = 2.47  # T4
$ws.Cells.Item(4, 21).Value = 1.7  # U4
$ws.Cells.Item(4, 22).Value = 1.93  # V4
$ws.Cells.Item(4, 23).Value = 9  # W4
$ws.Cells.Item(4, 24).Value = 15  # X4
$ws.Cells.Item(4, 25).Value = 10  # Y4
$ws.Cells.Item(4, 26).Value = 35  # Z4
$ws.Cells.Item(4, 27).Value = 23  # AA4
$ws.Cells.Item(4, 28).Value = 30  # AB4
$ws.Cells.Item(4, 29).Value = 9  # AC4
$ws.Cells.Item(4, 30).Value = 6  # AD4
$ws.Cells.Item(4, 31).Value = 13  # AE4
$ws.Cells.Item(4, 32).Value = 60  # AF4
$ws.Cells.Item(4, 33).Value = 450  # AG4
$ws.Cells.Item(4, 34).Value = 7.8  # AH4
$ws.Cells.Item(4, 35).Value = 12  # AI4
$ws.Cells.Item(4, 36).Value = 9.5  # AJ4
$ws.Cells.Item(4, 37).Value = 27  # AK4
$ws.Cells.Item(4, 38).Value = 21  # AL4
$ws.Cells.Item(4, 39).Value = 30  # AM4
$ws.Cells.Item(4, 40).Value = 4.75  # AN4
$ws.Cells.Item(4, 41).Value = 14.5  # AO4
$ws.Cells.Item(4, 42).Value = 19.5  # AP4
$ws.Cells.Item(4, 43).Value = 60  # AQ4
$ws.Cells.Item(4, 44).Value = 80  # AR4
$ws.Cells.Item(4, 45).Value = 200  # AS4
$ws.Cells.Item(4, 46).Value = 2.57  # AT4
$ws.Cells.Item(4, 47).Value = 6.5  # AU4
$ws.Cells.Item(4, 48).Value = 50  # AV4
$ws.Cells.Item(4, 49).Value = 4.45  # AW4
$ws.Cells.Item(4, 50).Value = 13  # AX4
$ws.Cells.Item(4, 51).Value = 20  # AY4
$ws.Cells.Item(4, 52).Value = 55  # AZ4
$ws.Cells.Item(4, 53).Value = 90  # BA4
$ws.Cells.Item(4, 54).Value = 250  # BB4
$ws.Cells.Item(4, 55).Value = 51  # BC4
$ws.Cells.Item(4, 56).Value = 51  # BD4
